$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the text values in B5:B8 (leave style intact, cells become blank)
$ws.Range("B5").Value = $null
$ws.Range("B6").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("B8").Value = $null
